$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-09-29 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-30 Saturday", 2) | Out-Null
$d.Content.Find.Execute("91-62=29", $true, $false, $false, $false, $false, $true, 1, $false, "24+55=79", 2) | Out-Null
$d.Content.Find.Execute("77-11=66", $true, $false, $false, $false, $false, $true, 1, $false, "5+41=46", 2) | Out-Null
$d.Content.Find.Execute("69-33=36", $true, $false, $false, $false, $false, $true, 1, $false, "86-6=80", 2) | Out-Null
$d.Content.Find.Execute("74-31=43", $true, $false, $false, $false, $false, $true, 1, $false, "37+11=48", 2) | Out-Null
$d.Content.Find.Execute("1+43=44", $true, $false, $false, $false, $false, $true, 1, $false, "98-83=15", 2) | Out-Null
$d.Content.Find.Execute("39+0=39", $true, $false, $false, $false, $false, $true, 1, $false, "85+11=96", 2) | Out-Null
$d.Content.Find.Execute("50+23=73", $true, $false, $false, $false, $false, $true, 1, $false, "15+40=55", 2) | Out-Null
$d.Content.Find.Execute("78-76=2", $true, $false, $false, $false, $false, $true, 1, $false, "60+25=85", 2) | Out-Null
$d.Content.Find.Execute("14+69=83", $true, $false, $false, $false, $false, $true, 1, $false, "85-80=5", 2) | Out-Null
$d.Content.Find.Execute("0+35=35", $true, $false, $false, $false, $false, $true, 1, $false, "21-4=17", 2) | Out-Null
$d.Content.Find.Execute("88-70=18", $true, $false, $false, $false, $false, $true, 1, $false, "33+19=52", 2) | Out-Null
$d.Content.Find.Execute("73-26=47", $true, $false, $false, $false, $false, $true, 1, $false, "27-25=2", 2) | Out-Null
$d.Content.Find.Execute("68-4=64", $true, $false, $false, $false, $false, $true, 1, $false, "73-55=18", 2) | Out-Null
$d.Content.Find.Execute("76-36=40", $true, $false, $false, $false, $false, $true, 1, $false, "87+1=88", 2) | Out-Null
$d.Content.Find.Execute("40+33=73", $true, $false, $false, $false, $false, $true, 1, $false, "79-51=28", 2) | Out-Null
$d.Content.Find.Execute("42+18=60", $true, $false, $false, $false, $false, $true, 1, $false, "60+37=97", 2) | Out-Null
$d.Content.Find.Execute("26+13=39", $true, $false, $false, $false, $false, $true, 1, $false, "23+75=98", 2) | Out-Null
$d.Content.Find.Execute("14+5=19", $true, $false, $false, $false, $false, $true, 1, $false, "47-33=14", 2) | Out-Null
$d.Content.Find.Execute("25+15=40", $true, $false, $false, $false, $false, $true, 1, $false, "52+0=52", 2) | Out-Null
$d.Content.Find.Execute("93-40=53", $true, $false, $false, $false, $false, $true, 1, $false, "43-31=12", 2) | Out-Null
$d.Content.Find.Execute("83-62=21", $true, $false, $false, $false, $false, $true, 1, $false, "21+10=31", 2) | Out-Null
$d.Content.Find.Execute("28+12=40", $true, $false, $false, $false, $false, $true, 1, $false, "78-47=31", 2) | Out-Null
$d.Content.Find.Execute("64-23=41", $true, $false, $false, $false, $false, $true, 1, $false, "27+67=94", 2) | Out-Null
$d.Content.Find.Execute("2+63=65", $true, $false, $false, $false, $false, $true, 1, $false, "32-12=20", 2) | Out-Null
$d.Content.Find.Execute("56+10=66", $true, $false, $false, $false, $false, $true, 1, $false, "65+18=83", 2) | Out-Null
$d.Content.Find.Execute("57-40=17", $true, $false, $false, $false, $false, $true, 1, $false, "5+80=85", 2) | Out-Null
$d.Content.Find.Execute("69+7=76", $true, $false, $false, $false, $false, $true, 1, $false, "25+52=77", 2) | Out-Null
$d.Content.Find.Execute("45+26=71", $true, $false, $false, $false, $false, $true, 1, $false, "65+16=81", 2) | Out-Null
$d.Content.Find.Execute("12-8=4", $true, $false, $false, $false, $false, $true, 1, $false, "0+97=97", 2) | Out-Null
$d.Content.Find.Execute("48+36=84", $true, $false, $false, $false, $false, $true, 1, $false, "29-23=6", 2) | Out-Null
$d.Content.Find.Execute("50+12=62", $true, $false, $false, $false, $false, $true, 1, $false, "85-19=66", 2) | Out-Null
$d.Content.Find.Execute("46+4=50", $true, $false, $false, $false, $false, $true, 1, $false, "13+51=64", 2) | Out-Null
$d.Content.Find.Execute("77+9=86", $true, $false, $false, $false, $false, $true, 1, $false, "52-41=11", 2) | Out-Null
$d.Content.Find.Execute("31+20=51", $true, $false, $false, $false, $false, $true, 1, $false, "29+29=58", 2) | Out-Null
$d.Content.Find.Execute("14+52=66", $true, $false, $false, $false, $false, $true, 1, $false, "77-13=64", 2) | Out-Null
$d.Content.Find.Execute("93-59=34", $true, $false, $false, $false, $false, $true, 1, $false, "97-71=26", 2) | Out-Null
$d.Content.Find.Execute("12+45=57", $true, $false, $false, $false, $false, $true, 1, $false, "40+7=47", 2) | Out-Null
$d.Content.Find.Execute("33+47=80", $true, $false, $false, $false, $false, $true, 1, $false, "54-7=47", 2) | Out-Null
$d.Content.Find.Execute("53-52=1", $true, $false, $false, $false, $false, $true, 1, $false, "93-7=86", 2) | Out-Null
$d.Content.Find.Execute("41-9=32", $true, $false, $false, $false, $false, $true, 1, $false, "39-3=36", 2) | Out-Null
$d.Content.Find.Execute("17-17=0", $true, $false, $false, $false, $false, $true, 1, $false, "63-55=8", 2) | Out-Null
$d.Content.Find.Execute("76-40=36", $true, $false, $false, $false, $false, $true, 1, $false, "77-71=6", 2) | Out-Null
$d.Content.Find.Execute("18+65=83", $true, $false, $false, $false, $false, $true, 1, $false, "93-32=61", 2) | Out-Null
$d.Content.Find.Execute("93+4=97", $true, $false, $false, $false, $false, $true, 1, $false, "36-22=14", 2) | Out-Null
$d.Content.Find.Execute("3+29=32", $true, $false, $false, $false, $false, $true, 1, $false, "16+37=53", 2) | Out-Null
$d.Content.Find.Execute("5+88=93", $true, $false, $false, $false, $false, $true, 1, $false, "17+78=95", 2) | Out-Null
$d.Content.Find.Execute("93-24=69", $true, $false, $false, $false, $false, $true, 1, $false, "73+17=90", 2) | Out-Null
$d.Content.Find.Execute("41-28=13", $true, $false, $false, $false, $false, $true, 1, $false, "30+26=56", 2) | Out-Null
$d.Content.Find.Execute("47-31=16", $true, $false, $false, $false, $false, $true, 1, $false, "43-28=15", 2) | Out-Null
$d.Content.Find.Execute("86-4=82", $true, $false, $false, $false, $false, $true, 1, $false, "12+25=37", 2) | Out-Null
$d.Content.Find.Execute("31-23=8", $true, $false, $false, $false, $false, $true, 1, $false, "72-14=58", 2) | Out-Null
$d.Content.Find.Execute("6+37=43", $true, $false, $false, $false, $false, $true, 1, $false, "73+17=90", 2) | Out-Null
$d.Content.Find.Execute("90-43=47", $true, $false, $false, $false, $false, $true, 1, $false, "8+16=24", 2) | Out-Null
$d.Content.Find.Execute("10-6=4", $true, $false, $false, $false, $false, $true, 1, $false, "4+46=50", 2) | Out-Null
$d.Content.Find.Execute("3+66=69", $true, $false, $false, $false, $false, $true, 1, $false, "60+30=90", 2) | Out-Null
$d.Content.Find.Execute("27+7=34", $true, $false, $false, $false, $false, $true, 1, $false, "5+92=97", 2) | Out-Null
$d.Content.Find.Execute("36+45=81", $true, $false, $false, $false, $false, $true, 1, $false, "31-25=6", 2) | Out-Null
$d.Content.Find.Execute("13+40=53", $true, $false, $false, $false, $false, $true, 1, $false, "2+51=53", 2) | Out-Null
$d.Content.Find.Execute("19-6=13", $true, $false, $false, $false, $false, $true, 1, $false, "8+23=31", 2) | Out-Null
$d.Content.Find.Execute("22+48=70", $true, $false, $false, $false, $false, $true, 1, $false, "71-68=3", 2) | Out-Null
$d.Content.Find.Execute("51+48=99", $true, $false, $false, $false, $false, $true, 1, $false, "15+62=77", 2) | Out-Null
$d.Content.Find.Execute("14+1=15", $true, $false, $false, $false, $false, $true, 1, $false, "58-43=15", 2) | Out-Null
$d.Content.Find.Execute("4+76=80", $true, $false, $false, $false, $false, $true, 1, $false, "16+79=95", 2) | Out-Null
$d.Content.Find.Execute("63-3=60", $true, $false, $false, $false, $false, $true, 1, $false, "2+56=58", 2) | Out-Null
$d.Content.Find.Execute("94-50=44", $true, $false, $false, $false, $false, $true, 1, $false, "6+28=34", 2) | Out-Null
$d.Content.Find.Execute("42-18=24", $true, $false, $false, $false, $false, $true, 1, $false, "65-15=50", 2) | Out-Null
$d.Content.Find.Execute("12+31=43", $true, $false, $false, $false, $false, $true, 1, $false, "4+47=51", 2) | Out-Null
$d.Content.Find.Execute("3+27=30", $true, $false, $false, $false, $false, $true, 1, $false, "0+87=87", 2) | Out-Null
$d.Content.Find.Execute("81-4=77", $true, $false, $false, $false, $false, $true, 1, $false, "55-15=40", 2) | Out-Null
$d.Content.Find.Execute("20+48=68", $true, $false, $false, $false, $false, $true, 1, $false, "97-92=5", 2) | Out-Null
$d.Content.Find.Execute("68+5=73", $true, $false, $false, $false, $false, $true, 1, $false, "35+27=62", 2) | Out-Null
$d.Content.Find.Execute("81+17=98", $true, $false, $false, $false, $false, $true, 1, $false, "72-43=29", 2) | Out-Null
$d.Content.Find.Execute("71+10=81", $true, $false, $false, $false, $false, $true, 1, $false, "42+49=91", 2) | Out-Null
$d.Content.Find.Execute("7+89=96", $true, $false, $false, $false, $false, $true, 1, $false, "14+43=57", 2) | Out-Null
$d.Content.Find.Execute("41+27=68", $true, $false, $false, $false, $false, $true, 1, $false, "99-55=44", 2) | Out-Null
$d.Content.Find.Execute("17+33=50", $true, $false, $false, $false, $false, $true, 1, $false, "77+10=87", 2) | Out-Null
$d.Content.Find.Execute("8+8=16", $true, $false, $false, $false, $false, $true, 1, $false, "38+49=87", 2) | Out-Null
$d.Content.Find.Execute("5+51=56", $true, $false, $false, $false, $false, $true, 1, $false, "3+17=20", 2) | Out-Null
$d.Content.Find.Execute("23-14=9", $true, $false, $false, $false, $false, $true, 1, $false, "49+23=72", 2) | Out-Null
$d.Content.Find.Execute("66-58=8", $true, $false, $false, $false, $false, $true, 1, $false, "24-16=8", 2) | Out-Null
$d.Content.Find.Execute("18+28=46", $true, $false, $false, $false, $false, $true, 1, $false, "94-56=38", 2) | Out-Null
$d.Content.Find.Execute("53+3=56", $true, $false, $false, $false, $false, $true, 1, $false, "18+50=68", 2) | Out-Null
$d.Content.Find.Execute("83-13=70", $true, $false, $false, $false, $false, $true, 1, $false, "13+5=18", 2) | Out-Null
$d.Content.Find.Execute("27+46=73", $true, $false, $false, $false, $false, $true, 1, $false, "11+29=40", 2) | Out-Null
$d.Content.Find.Execute("18+60=78", $true, $false, $false, $false, $false, $true, 1, $false, "0+86=86", 2) | Out-Null
$d.Content.Find.Execute("17+63=80", $true, $false, $false, $false, $false, $true, 1, $false, "41-36=5", 2) | Out-Null
$d.Content.Find.Execute("60+22=82", $true, $false, $false, $false, $false, $true, 1, $false, "49-12=37", 2) | Out-Null
$d.Content.Find.Execute("79-39=40", $true, $false, $false, $false, $false, $true, 1, $false, "92-66=26", 2) | Out-Null
$d.Content.Find.Execute("37+34=71", $true, $false, $false, $false, $false, $true, 1, $false, "42-0=42", 2) | Out-Null
$d.Content.Find.Execute("21-20=1", $true, $false, $false, $false, $false, $true, 1, $false, "45-32=13", 2) | Out-Null
$d.Content.Find.Execute("9+85=94", $true, $false, $false, $false, $false, $true, 1, $false, "31+63=94", 2) | Out-Null
$d.Content.Find.Execute("95-31=64", $true, $false, $false, $false, $false, $true, 1, $false, "52-32=20", 2) | Out-Null
$d.Content.Find.Execute("24-5=19", $true, $false, $false, $false, $false, $true, 1, $false, "96-64=32", 2) | Out-Null
$d.Content.Find.Execute("87-39=48", $true, $false, $false, $false, $false, $true, 1, $false, "11+88=99", 2) | Out-Null
$d.Content.Find.Execute("85-41=44", $true, $false, $false, $false, $false, $true, 1, $false, "69-6=63", 2) | Out-Null
$d.Content.Find.Execute("27+45=72", $true, $false, $false, $false, $false, $true, 1, $false, "76+19=95", 2) | Out-Null
$d.Content.Find.Execute("91-48=43", $true, $false, $false, $false, $false, $true, 1, $false, "86-36=50", 2) | Out-Null
$d.Content.Find.Execute("85-30=55", $true, $false, $false, $false, $false, $true, 1, $false, "46+19=65", 2) | Out-Null
$d.Content.Find.Execute("85+0=85", $true, $false, $false, $false, $false, $true, 1, $false, "94-91=3", 2) | Out-Null
$d.Content.Find.Execute("25+70=95", $true, $false, $false, $false, $false, $true, 1, $false, "18+49=67", 2) | Out-Null
